# Auto-generated Excel COM-interop script applying the OOXML diff
# to Sheets/Malboro_Profits.xlsx (workbook sheets ALC, ARM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 308.44446
$ws.Cells.Item(2, 9).Value = 130.66667
$ws.Cells.Item(2, 10).Value = 664
$ws.Cells.Item(2, 11).Value = 130.66667
$ws.Cells.Item(2, 12).Value = 664
$ws.Cells.Item(2, 13).Value = -17.66667000000001
$ws.Cells.Item(2, 14).Value = -890
$ws.Cells.Item(40, 8).Value = 1005.4286
$ws.Cells.Item(40, 10).Value = 1024.5
$ws.Cells.Item(40, 12).Value = 1024.5
$ws.Cells.Item(40, 14).Value = -1374.5
$ws.Cells.Item(62, 8).Value = 4997.2
$ws.Cells.Item(62, 10).Value = 5000
$ws.Cells.Item(62, 12).Value = 5000
$ws.Cells.Item(62, 14).Value = -6248
$ws.Cells.Item(65, 8).Value = 4997.2
$ws.Cells.Item(65, 10).Value = 5000
$ws.Cells.Item(65, 12).Value = 25000
$ws.Cells.Item(65, 14).Value = -31240
$ws.Cells.Item(132, 8).Value = 7275.8887
$ws.Cells.Item(132, 9).Value = 8450.666999999999
$ws.Cells.Item(132, 10).Value = 1402
$ws.Cells.Item(132, 11).Value = 25352.001
$ws.Cells.Item(132, 12).Value = 4206
$ws.Cells.Item(132, 13).Value = -22822.001
$ws.Cells.Item(132, 14).Value = -9266
$ws.Cells.Item(137, 8).Value = 9889.559999999999
$ws.Cells.Item(137, 9).Value = 1395.8572
$ws.Cells.Item(137, 10).Value = 20699.727
$ws.Cells.Item(137, 11).Value = 4187.571599999999
$ws.Cells.Item(137, 12).Value = 62099.181
$ws.Cells.Item(137, 13).Value = -1637.571599999999
$ws.Cells.Item(137, 14).Value = -67199.181
$ws.Cells.Item(141, 8).Value = 4511.2173
$ws.Cells.Item(141, 9).Value = 4649.5884
$ws.Cells.Item(141, 11).Value = 13948.7652
$ws.Cells.Item(141, 13).Value = -8768.765199999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 573756.25
$ws.Cells.Item(61, 9).Value = 2462.762
$ws.Cells.Item(61, 11).Value = 2462.762
$ws.Cells.Item(61, 13).Value = -2250.762
$ws.Cells.Item(132, 8).Value = 6625416
$ws.Cells.Item(132, 9).Value = 3015.5334
$ws.Cells.Item(132, 11).Value = 9046.600199999999
$ws.Cells.Item(132, 13).Value = -6516.600199999999
$ws.Cells.Item(136, 8).Value = 573756.25
$ws.Cells.Item(136, 9).Value = 2462.762
$ws.Cells.Item(136, 11).Value = 7388.286
$ws.Cells.Item(136, 13).Value = -4838.286

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(12, 8).Value = 417608.9
$ws.Cells.Item(12, 9).Value = 1027.909
$ws.Cells.Item(12, 10).Value = 5000000
$ws.Cells.Item(12, 11).Value = 1027.909
$ws.Cells.Item(12, 12).Value = 5000000
$ws.Cells.Item(12, 13).Value = -857.9090000000001
$ws.Cells.Item(12, 14).Value = -5000340
$ws.Cells.Item(31, 8).Value = 9967.886
$ws.Cells.Item(31, 9).Value = 1446.0344
$ws.Cells.Item(31, 11).Value = 1446.0344
$ws.Cells.Item(31, 13).Value = -1151.0344
$ws.Cells.Item(34, 8).Value = 9967.886
$ws.Cells.Item(34, 9).Value = 1446.0344
$ws.Cells.Item(34, 11).Value = 1446.0344
$ws.Cells.Item(34, 13).Value = -1244.0344
$ws.Cells.Item(54, 8).Value = 25000
$ws.Cells.Item(54, 9).Value = 25000
$ws.Cells.Item(54, 11).Value = 25000
$ws.Cells.Item(54, 13).Value = -24342
$ws.Cells.Item(58, 8).Value = 21282
$ws.Cells.Item(58, 9).Value = 11485.125
$ws.Cells.Item(58, 11).Value = 11485.125
$ws.Cells.Item(58, 13).Value = -11282.125
$ws.Cells.Item(62, 8).Value = 6191.1665
$ws.Cells.Item(62, 9).Value = 6185.2856
$ws.Cells.Item(62, 11).Value = 6185.2856
$ws.Cells.Item(62, 13).Value = -5561.2856
$ws.Cells.Item(65, 8).Value = 6191.1665
$ws.Cells.Item(65, 9).Value = 6185.2856
$ws.Cells.Item(65, 11).Value = 30926.428
$ws.Cells.Item(65, 13).Value = -27806.428
$ws.Cells.Item(86, 8).Value = 19935.75
$ws.Cells.Item(86, 9).Value = 23248.166
$ws.Cells.Item(86, 11).Value = 23248.166
$ws.Cells.Item(86, 13).Value = -22125.166
$ws.Cells.Item(89, 8).Value = 19935.75
$ws.Cells.Item(89, 9).Value = 23248.166
$ws.Cells.Item(89, 11).Value = 116240.83
$ws.Cells.Item(89, 13).Value = -110624.83
$ws.Cells.Item(132, 8).Value = 51473348
$ws.Cells.Item(132, 9).Value = 3086.625
$ws.Cells.Item(132, 11).Value = 9259.875
$ws.Cells.Item(132, 13).Value = -6729.875
$ws.Cells.Item(134, 8).Value = 28577402
$ws.Cells.Item(134, 9).Value = 2401.318
$ws.Cells.Item(134, 11).Value = 7203.954000000001
$ws.Cells.Item(134, 13).Value = -4668.954000000001
$ws.Cells.Item(136, 8).Value = 21282
$ws.Cells.Item(136, 9).Value = 11485.125
$ws.Cells.Item(136, 11).Value = 34455.375
$ws.Cells.Item(136, 13).Value = -31905.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 1877.6666
$ws.Cells.Item(107, 9).Value = 1175
$ws.Cells.Item(107, 10).Value = 2439.8
$ws.Cells.Item(107, 11).Value = 3525
$ws.Cells.Item(107, 12).Value = 7319.400000000001
$ws.Cells.Item(107, 13).Value = -1605
$ws.Cells.Item(107, 14).Value = -11159.4
$ws.Cells.Item(131, 8).Value = 10877.747
$ws.Cells.Item(131, 9).Value = 9913.799999999999
$ws.Cells.Item(131, 10).Value = 10991.153
$ws.Cells.Item(131, 11).Value = 29741.4
$ws.Cells.Item(131, 12).Value = 32973.459
$ws.Cells.Item(131, 13).Value = -24701.4
$ws.Cells.Item(131, 14).Value = -43053.459
$ws.Cells.Item(133, 8).Value = 5015
$ws.Cells.Item(133, 9).Value = 3686.6667
$ws.Cells.Item(133, 11).Value = 11060.0001
$ws.Cells.Item(133, 13).Value = -6000.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(46, 8).Value = 6400
$ws.Cells.Item(46, 10).Value = 0
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(69, 8).Value = 43282.715
$ws.Cells.Item(69, 10).Value = 43282.715
$ws.Cells.Item(69, 12).Value = 43282.715
$ws.Cells.Item(69, 14).Value = -44780.715
$ws.Cells.Item(72, 8).Value = 43282.715
$ws.Cells.Item(72, 10).Value = 43282.715
$ws.Cells.Item(72, 12).Value = 129848.145
$ws.Cells.Item(72, 14).Value = -137336.145
$ws.Cells.Item(102, 8).Value = 5744.6294
$ws.Cells.Item(102, 9).Value = 6456.4287
$ws.Cells.Item(102, 11).Value = 6456.4287
$ws.Cells.Item(102, 13).Value = -4834.4287
$ws.Cells.Item(126, 8).Value = 9610.272000000001
$ws.Cells.Item(126, 9).Value = 11550
$ws.Cells.Item(126, 11).Value = 34650
$ws.Cells.Item(126, 13).Value = -32180
$ws.Cells.Item(132, 8).Value = 2445723.8
$ws.Cells.Item(132, 9).Value = 6950
$ws.Cells.Item(132, 11).Value = 20850
$ws.Cells.Item(132, 13).Value = -18320
$ws.Cells.Item(135, 8).Value = 79000
$ws.Cells.Item(135, 10).Value = 79000
$ws.Cells.Item(135, 12).Value = 79000
$ws.Cells.Item(135, 14).Value = -89140
$ws.Cells.Item(46, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2500
$ws.Cells.Item(40, 10).Value = 2500
$ws.Cells.Item(40, 12).Value = 2500
$ws.Cells.Item(40, 14).Value = -2772
$ws.Cells.Item(44, 8).Value = 15250
$ws.Cells.Item(44, 10).Value = 15250
$ws.Cells.Item(44, 12).Value = 15250
$ws.Cells.Item(44, 14).Value = -16162
$ws.Cells.Item(46, 8).Value = 6833
$ws.Cells.Item(46, 10).Value = 9999.5
$ws.Cells.Item(46, 12).Value = 9999.5
$ws.Cells.Item(46, 14).Value = -10375.5
$ws.Cells.Item(55, 8).Value = 2054.15
$ws.Cells.Item(55, 9).Value = 2356.2856
$ws.Cells.Item(55, 10).Value = 1891.4615
$ws.Cells.Item(55, 11).Value = 2356.2856
$ws.Cells.Item(55, 12).Value = 1891.4615
$ws.Cells.Item(55, 13).Value = -2183.2856
$ws.Cells.Item(55, 14).Value = -2237.4615
$ws.Cells.Item(136, 8).Value = 1147288.9
$ws.Cells.Item(136, 10).Value = 2657654.5
$ws.Cells.Item(136, 12).Value = 7972963.5
$ws.Cells.Item(136, 14).Value = -7978063.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(49, 8).Value = 29999.857
$ws.Cells.Item(49, 9).Value = 29999.857
$ws.Cells.Item(49, 11).Value = 29999.857
$ws.Cells.Item(49, 13).Value = -29769.857
$ws.Cells.Item(70, 8).Value = 7800
$ws.Cells.Item(70, 10).Value = 7800
$ws.Cells.Item(70, 12).Value = 7800
$ws.Cells.Item(70, 14).Value = -8430
$ws.Cells.Item(73, 8).Value = 7800
$ws.Cells.Item(73, 10).Value = 7800
$ws.Cells.Item(73, 12).Value = 7800
$ws.Cells.Item(73, 14).Value = -9984
$ws.Cells.Item(132, 8).Value = 1016382.8
$ws.Cells.Item(132, 9).Value = 15235.143
$ws.Cells.Item(132, 11).Value = 45705.429
$ws.Cells.Item(132, 13).Value = -43175.429
$ws.Cells.Item(136, 8).Value = 521961.06
$ws.Cells.Item(136, 9).Value = 2213.0667
$ws.Cells.Item(136, 11).Value = 6639.2001
$ws.Cells.Item(136, 13).Value = -4089.2001
